$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape
$ws.Range('D2').Value = '24.674.42'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '1.688.97'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3937'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.99%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4051'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.003'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.57'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08834'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.236'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.021'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001317'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '1.689.92'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07014'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.988'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.006'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.48%  '
$ws.Range('D24').Value = '24.657.64'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.291'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +8.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.355'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('E28').Value = '  +1.62%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.179'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '135.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.622'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.87%  '
$ws.Range('D32').Value = '1.875.80'
$ws.Range('E32').Value = '  -0.50%  '
$ws.Range('E33').Value = '  -2.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08527'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.106'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.23%  '
$ws.Range('E36').Value = '  +1.74%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2728'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.885'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.63%  '
$ws.Range('E39').Value = '  -3.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09195'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02714'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.462'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7624'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.588'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.7132'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.215'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('E48').Value = '  +0.44%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.319'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '139.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07966'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.44%  '
